$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "Student_ID"
$ws.Range("B1").Value = "First_Name"
$ws.Range("C1").Value = "Middle_Initial"
$ws.Range("D1").Value = "Last_Name"
$ws.Range("E1").Value = "Date_of_Birth"
$ws.Range("F1").Value = "CAIR_Program_Name"
$ws.Range("G1").Value = "Hours of Participation"
$ws.Range("H1").Value = "Notes"

# Apply font formatting to A1:E1 (Calibri 12)
$rng = $ws.Range("A1:E1")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 12

# Row height
$ws.Rows.Item(1).RowHeight = 15.75

# Selection
$ws.Range("K9").Select()
